# Updates the Price (column D) and Volume(1h) (column E) cells of the
# cryptocurrency listing on the active worksheet to reflect the latest
# scraped values, matching a GitHub Actions scheduled refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values such as "1.007" or "219.40" look like plain numbers to Excel and
# would silently be converted from text to a numeric value if assigned
# directly. The source data keeps them as text (note some prices such as
# "26.218.95" use multiple "." separators and are never numeric). To keep
# every updated cell as text - exactly like the rest of the column - values
# that would be auto-recognised as numbers are entered with a leading
# apostrophe (Excel's standard way to force text entry) and the cell style
# is then reset back to "Normal" so no visible formatting change remains.
$numericPattern = '^\s*[+-]?[0-9]+(\.[0-9]+)?\s*$'

function Set-TextValue($cell, $value) {
    if ($value -match $numericPattern) {
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

Set-TextValue $ws.Cells.Item(2, 4) "26.218.95"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.22%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.660.67"
Set-TextValue $ws.Cells.Item(3, 5) "  -0.18%  "
Set-TextValue $ws.Cells.Item(4, 4) "1.006"
Set-TextValue $ws.Cells.Item(4, 5) "  -0.43%  "
Set-TextValue $ws.Cells.Item(5, 4) "219.40"
Set-TextValue $ws.Cells.Item(5, 5) "  +0.28%  "
Set-TextValue $ws.Cells.Item(6, 4) "0.5260"
Set-TextValue $ws.Cells.Item(6, 5) "  -0.90%  "
Set-TextValue $ws.Cells.Item(7, 5) "  -0.40%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.2636"
Set-TextValue $ws.Cells.Item(8, 5) "  +0.03%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.06320"
Set-TextValue $ws.Cells.Item(9, 5) "  -0.67%  "
Set-TextValue $ws.Cells.Item(10, 4) "20.65"
Set-TextValue $ws.Cells.Item(10, 5) "  +0.51%  "
Set-TextValue $ws.Cells.Item(11, 5) "  -0.53%  "
Set-TextValue $ws.Cells.Item(12, 4) "4.501"
Set-TextValue $ws.Cells.Item(12, 5) "  -1.35%  "
Set-TextValue $ws.Cells.Item(13, 4) "1.664.90"
Set-TextValue $ws.Cells.Item(13, 5) "  +0.07%  "
Set-TextValue $ws.Cells.Item(14, 4) "1.888.89"
Set-TextValue $ws.Cells.Item(14, 5) "  -0.16%  "
Set-TextValue $ws.Cells.Item(15, 4) "0.5558"
Set-TextValue $ws.Cells.Item(15, 5) "  +0.56%  "
Set-TextValue $ws.Cells.Item(16, 4) "0.0₅8029"
Set-TextValue $ws.Cells.Item(16, 5) "  -1.81%  "
Set-TextValue $ws.Cells.Item(17, 4) "65.26"
Set-TextValue $ws.Cells.Item(17, 5) "  -0.55%  "
Set-TextValue $ws.Cells.Item(18, 4) "26.221.72"
Set-TextValue $ws.Cells.Item(19, 4) "1.006"
Set-TextValue $ws.Cells.Item(19, 5) "  -0.38%  "
Set-TextValue $ws.Cells.Item(20, 4) "4.644"
Set-TextValue $ws.Cells.Item(20, 5) "  -0.66%  "
Set-TextValue $ws.Cells.Item(21, 4) "196.78"
Set-TextValue $ws.Cells.Item(21, 5) "  +2.07%  "
Set-TextValue $ws.Cells.Item(22, 4) "10.15"
Set-TextValue $ws.Cells.Item(22, 5) "  -0.69%  "
Set-TextValue $ws.Cells.Item(23, 4) "5.973"
Set-TextValue $ws.Cells.Item(23, 5) "  -0.97%  "
Set-TextValue $ws.Cells.Item(24, 4) "1.007"
Set-TextValue $ws.Cells.Item(24, 5) "  -0.44%  "
Set-TextValue $ws.Cells.Item(25, 4) "145.78"
Set-TextValue $ws.Cells.Item(25, 5) "  +0.99%  "
Set-TextValue $ws.Cells.Item(26, 4) "0.1208"
Set-TextValue $ws.Cells.Item(26, 5) "  -1.49%  "
Set-TextValue $ws.Cells.Item(27, 4) "7.159"
Set-TextValue $ws.Cells.Item(28, 4) "16.05"
Set-TextValue $ws.Cells.Item(28, 5) "  -0.13%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.509"
Set-TextValue $ws.Cells.Item(29, 5) "  +2.13%  "
Set-TextValue $ws.Cells.Item(30, 5) "  -2.41%  "
Set-TextValue $ws.Cells.Item(31, 5) "  -0.24%  "
Set-TextValue $ws.Cells.Item(32, 4) "3.491"
Set-TextValue $ws.Cells.Item(32, 5) "  -2.88%  "
Set-TextValue $ws.Cells.Item(33, 4) "3.355"
Set-TextValue $ws.Cells.Item(33, 5) "  +2.27%  "
Set-TextValue $ws.Cells.Item(34, 4) "1.585"
Set-TextValue $ws.Cells.Item(34, 5) "  -1.84%  "
Set-TextValue $ws.Cells.Item(35, 5) "  -0.59%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.9544"
Set-TextValue $ws.Cells.Item(37, 4) "2.420"
Set-TextValue $ws.Cells.Item(37, 5) "  -0.21%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.5735"
Set-TextValue $ws.Cells.Item(38, 5) "  -1.03%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.01594"
Set-TextValue $ws.Cells.Item(39, 5) "  -0.53%  "
Set-TextValue $ws.Cells.Item(40, 4) "5.964"
Set-TextValue $ws.Cells.Item(40, 5) "  +1.76%  "
Set-TextValue $ws.Cells.Item(41, 4) "1.061.14"
Set-TextValue $ws.Cells.Item(41, 5) "  +1.38%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.8517"
Set-TextValue $ws.Cells.Item(42, 5) "  -1.65%  "
Set-TextValue $ws.Cells.Item(43, 5) "  -0.37%  "
Set-TextValue $ws.Cells.Item(44, 4) "103.17"
Set-TextValue $ws.Cells.Item(45, 4) "1.799.86"
Set-TextValue $ws.Cells.Item(45, 5) "  -0.27%  "
Set-TextValue $ws.Cells.Item(46, 4) "58.28"
Set-TextValue $ws.Cells.Item(46, 5) "  +1.53%  "
Set-TextValue $ws.Cells.Item(47, 4) "1.010"
Set-TextValue $ws.Cells.Item(47, 5) "  -0.27%  "
Set-TextValue $ws.Cells.Item(48, 4) "0.4408"
Set-TextValue $ws.Cells.Item(48, 5) "  +0.58%  "
Set-TextValue $ws.Cells.Item(49, 4) "8.005"
Set-TextValue $ws.Cells.Item(49, 5) "  +0.07%  "
Set-TextValue $ws.Cells.Item(50, 4) "0.05202"
Set-TextValue $ws.Cells.Item(50, 5) "  +0.69%  "
Set-TextValue $ws.Cells.Item(51, 5) "  -5.20%  "
